$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data updates (Client Id, User Name, Exam Password, First Name, Last Name, Candidate ID)
$ws.Range("A2").Value = "oNxCC145"
$ws.Range("C2").Value = "nrstfbf83"
$ws.Range("D2").Value = 'sg6S#$J9'
$ws.Range("F2").Value = "CaSXNnUx"
$ws.Range("G2").Value = "bMuI"
$ws.Range("B2").Value = 23082431
